$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 41203
$ws.Cells.Item(2, 2).Value = "Ana Clara Novaes"
$ws.Cells.Item(2, 3).Value = "Operações"
$ws.Cells.Item(2, 4).Value = "Consulta médica"
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 45095
$ws.Cells.Item(2, 7).Value = 11081.15

# Row 3
$ws.Cells.Item(3, 1).Value = 42647
$ws.Cells.Item(3, 2).Value = "Yasmin Barros"
$ws.Cells.Item(3, 3).Value = "TI"
$ws.Cells.Item(3, 4).Value = "Doença"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 45081
$ws.Cells.Item(3, 7).Value = 3882.83

# Row 4
$ws.Cells.Item(4, 1).Value = 35999
$ws.Cells.Item(4, 2).Value = "Sr. Davi Farias"
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 7).Value = 3540.8

# Row 5
$ws.Cells.Item(5, 1).Value = 14554
$ws.Cells.Item(5, 2).Value = "Ana Sophia Martins"
$ws.Cells.Item(5, 3).Value = "Jurídico"
$ws.Cells.Item(5, 4).Value = "Doença"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 45099
$ws.Cells.Item(5, 7).Value = 11738.56

# Row 6
$ws.Cells.Item(6, 1).Value = 53013
$ws.Cells.Item(6, 2).Value = "Benício Ramos"
$ws.Cells.Item(6, 4).Value = "Doença"
$ws.Cells.Item(6, 6).Value = 45087
$ws.Cells.Item(6, 7).Value = 2558.18

# Row 7
$ws.Cells.Item(7, 1).Value = 66510
$ws.Cells.Item(7, 2).Value = "Alice da Conceição"
$ws.Cells.Item(7, 3).Value = "Operações"
$ws.Cells.Item(7, 4).Value = "Viagem de negócios"
$ws.Cells.Item(7, 5).Value = 6
$ws.Cells.Item(7, 6).Value = 45087
$ws.Cells.Item(7, 7).Value = 4999.9

# Row 8
$ws.Cells.Item(8, 1).Value = 91107
$ws.Cells.Item(8, 2).Value = "Esther Rodrigues"
$ws.Cells.Item(8, 3).Value = "Operações"
$ws.Cells.Item(8, 4).Value = "Problemas pessoais"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 45078
$ws.Cells.Item(8, 7).Value = 4081.16

# Row 9
$ws.Cells.Item(9, 1).Value = 95915
$ws.Cells.Item(9, 2).Value = "Yago Barbosa"
$ws.Cells.Item(9, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(9, 4).Value = "Doença"
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 45101
$ws.Cells.Item(9, 7).Value = 10410.72

# Row 10
$ws.Cells.Item(10, 1).Value = 32422
$ws.Cells.Item(10, 2).Value = "Matheus da Rosa"
$ws.Cells.Item(10, 3).Value = "Operações"
$ws.Cells.Item(10, 4).Value = "Consulta médica"
$ws.Cells.Item(10, 5).Value = 6
$ws.Cells.Item(10, 6).Value = 45078
$ws.Cells.Item(10, 7).Value = 8157.49

# Row 11
$ws.Cells.Item(11, 1).Value = 17110
$ws.Cells.Item(11, 2).Value = "Benício Mendes"
$ws.Cells.Item(11, 3).Value = "TI"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 45085
$ws.Cells.Item(11, 7).Value = 10660.61
